# Weekly update: insert a new price record at the top of the Acelga
# history table (row 93), shifting existing rows 93-109 down to 94-110.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93 (pushes old row 93..109 down to 94..110,
# and copies the formatting - including the date style on column D - from
# the row above, matching Excel's native "Insert Row" behaviour).
$ws.Rows.Item(93).Insert()

# Populate the new row 93 with the latest weekly record.
$ws.Cells.Item(93, 1).Value = 1
$ws.Cells.Item(93, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(93, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(93, 4).Value = 45244
$ws.Cells.Item(93, 5).Value = 15
$ws.Cells.Item(93, 6).Value = 100112009
$ws.Cells.Item(93, 7).Value = "Acelga"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Segunda"
$ws.Cells.Item(93, 10).Value = 200
$ws.Cells.Item(93, 11).Value = 1000
$ws.Cells.Item(93, 12).Value = 1500
$ws.Cells.Item(93, 13).Value = 1250
$ws.Cells.Item(93, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(93, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(93, 16).Value = 417
$ws.Cells.Item(93, 17).Value = 3
$ws.Cells.Item(93, 18).Value = "Hortaliza"
